# 1.2.1.1a workbook update:
#  - Re-capitalise / rewrap the "by sex" / "by age" / "by territory" / "by
#    type of activity" section header labels (columns A/B/C) for rows
#    6,7,8,9,22,23,24,25.
#  - Add a new "2020" data column (I) with the reported values, copying the
#    number formatting from the corresponding 2019 column (H), except for
#    row 46 (the "…" placeholder row) which copies from the already-numeric
#    F column, and row 37 which gets a dedicated right-aligned style.
#  - Leave the selection on B28, matching the saved state in the source file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Section header label updates -----------------------------------
$ws.Range("A6").Value = "Жынысы боюнча"
$ws.Range("B6").Value = "По полу"
$ws.Range("C6").Value = "By sex"

$ws.Range("A7").Value = "Эркектер"
$ws.Range("B7").Value = "Мужчины"
$ws.Range("C7").Value = "Men"

$ws.Range("A8").Value = "Аялдар"
$ws.Range("B8").Value = "Женщины"
$ws.Range("C8").Value = "Women"

$ws.Range("A9").Value = "Жаш-курагы бонюча"
$ws.Range("B9").Value = "По возрастным группам:"
$ws.Range("C9").Value = "By age group"

$ws.Range("A22").Value = "Жерлер"
$ws.Range("B22").Value = "Местность"
$ws.Range("C22").Value = "Urbanization"

$ws.Range("A23").Value = "Шаар"
$ws.Range("B23").Value = "Город"
$ws.Range("C23").Value = "Urban"

$ws.Range("A24").Value = "Айыл"
$ws.Range("B24").Value = "Село"
$ws.Range("C24").Value = "Rural"

$ws.Range("A25").Value = "Экономиканын тармактары боюнча"
$ws.Range("B25").Value = "По отраслям экономики"
$ws.Range("C25").Value = "By type of  activity "

# ---- New "2020" column (I) -------------------------------------------
$values2020 = @{
  4  = 2020
  5  = 17.2
  7  = 18.7
  8  = 14.9
  10 = 11.4
  11 = 8.9
  12 = 15.3
  13 = 22.1
  14 = 27.5
  15 = 25.2
  16 = 16.1
  17 = 12.3
  18 = 8.9
  19 = 10.2
  20 = 9.4
  21 = 5.1
  23 = 12.4
  24 = 20.1
  26 = 21.1
  27 = 40.1
  28 = 13.5
  29 = 18.7
  30 = 24.4
  31 = 21.6
  32 = 10.1
  33 = 19.4
  34 = 11.6
  35 = 5.3
  36 = 20.5
  37 = 9
  38 = 21
  39 = 12
  40 = 13.1
  41 = 22
  42 = 22.9
  43 = 12.6
  44 = 12.5
  45 = 36
  46 = 10.7
}

for ($r = 4; $r -le 46; $r++) {
  $srcCol = "H"
  if ($r -eq 46) { $srcCol = "F" }
  [void]$ws.Range("$srcCol$r").Copy($ws.Range("I$r"))
  if ($values2020.ContainsKey($r)) {
    $ws.Range("I$r").Value = $values2020[$r]
  }
}

# Row 37 gets its own right-aligned number style (new cellXfs entry) rather
# than the plain style inherited from H37.
$ws.Range("I37").HorizontalAlignment = -4152

# ---- Restore saved selection ------------------------------------------
[void]$ws.Range("B28").Select()
